# Generate Report for Handoff
# Updates the "b.md" row (row 3) in the Overview, zh-cn and de-de sheets to
# reflect that the file is now "Ready for handoff" with a fresh handoff
# file/timestamp, instead of the previous "Handed back" status.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# Row 3 corresponds to b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 02:33:17"

# ---------------------------------------------------------------------
# Helper: update the per-language detail sheet (row 3 = b.md)
#   C -> Status
#   D -> Latest Handoff File (+ hyperlink display text, same target)
#   E -> Latest Handoff Datetime
# ---------------------------------------------------------------------
function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$HandoffFileName,
        [string]$HandoffDatetime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $HandoffFileName
    $ws.Range("E3").Value = $HandoffDatetime

    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq 3 -and $hl.Range.Column -eq 4) {
            $hl.TextToDisplay = $HandoffFileName
        }
    }
}

Update-LanguageSheet "zh-cn" "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" "2016-03-23 02:33:13"
Update-LanguageSheet "de-de" "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" "2016-03-23 02:33:17"
